$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells are stored as text in the source data (values like
# "97.505.65" or "242.39" are not valid numbers). Force NumberFormat to "@" (Text)
# before assigning so Excel keeps them as text instead of auto-coercing to a
# number (which would also silently drop meaningful trailing zeros, e.g. "5.50").

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.505.65"
$ws.Range("E2").Value = "  +4.68%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.125.71"
$ws.Range("E3").Value = "  +0.09%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.39"
$ws.Range("E5").Value = "  +2.08%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "610.05"
$ws.Range("E6").Value = "  -0.64%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.12"
$ws.Range("E7").Value = "  +1.91%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.382"
$ws.Range("E8").Value = "  -2.07%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.125.58"
$ws.Range("E10").Value = "  +0.17%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.775"
$ws.Range("E11").Value = "  -7.10%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.198"
$ws.Range("E12").Value = "  -0.22%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "97.137.09"
$ws.Range("E13").Value = "  +4.45%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000238"
$ws.Range("E14").Value = "  -2.89%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.50"
$ws.Range("E15").Value = "  +1.01%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "33.86"
$ws.Range("E16").Value = "  -4.30%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.709.77"
$ws.Range("E17").Value = "  +0.06%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.132.57"
$ws.Range("E18").Value = "  +0.80%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.49"
$ws.Range("E19").Value = "  -6.60%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "516.53"
$ws.Range("E20").Value = "  +16.66%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.42"
$ws.Range("E21").Value = "  -2.69%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.62"
$ws.Range("E22").Value = "  -6.97%  "

# Row 23
$ws.Range("E23").Value = "  -5.20%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.75"
$ws.Range("E24").Value = "  -4.35%  "

# Row 25
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "88.64"
$ws.Range("E25").Value = "  +2.95%  "

# Row 26
$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.46"
$ws.Range("E26").Value = "  -4.67%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.54"
$ws.Range("E27").Value = "  -12.39%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.298.68"
$ws.Range("E28").Value = "  +0.29%  "

# Row 29
$ws.Range("E29").Value = "  +0.30%  "

# Row 30
$ws.Range("E30").Value = "  +3.90%  "

# Row 31
$ws.Range("E31").Value = "  -4.39%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.123"
$ws.Range("E32").Value = "  -0.75%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.986"
$ws.Range("E33").Value = "  -1.91%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "8.95"
$ws.Range("E34").Value = "  -3.68%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.43"
$ws.Range("E35").Value = "  +1.36%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.151"
$ws.Range("E36").Value = "  -5.42%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.25"
$ws.Range("E37").Value = "  -9.76%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.87"
$ws.Range("E38").Value = "  -1.93%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.23"
$ws.Range("E39").Value = "  +0.98%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.433"
$ws.Range("E40").Value = "  -4.14%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "465.37"
$ws.Range("E41").Value = "  -2.40%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.57"
$ws.Range("E42").Value = "  -10.56%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.21"
$ws.Range("E43").Value = "  -6.82%  "

# Row 44
$ws.Range("E44").Value = "  -0.02%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.10"
$ws.Range("E45").Value = "  -6.98%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "162.01"
$ws.Range("E46").Value = "  +2.04%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.91"
$ws.Range("E47").Value = "  +2.75%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.693"
$ws.Range("E48").Value = "  -1.36%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.46"
$ws.Range("E49").Value = "  +0.37%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.24"
$ws.Range("E50").Value = "  +0.38%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  +0.01%  "
